$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1004.9
$ws.Range("I19").Value = 1009.6
$ws.Range("J19").Value = 1000.2
$ws.Range("K19").Value = 1009.6
$ws.Range("L19").Value = 1000.2
$ws.Range("M19").Value = -834.6
$ws.Range("N19").Value = -1350.2

$ws.Range("H51").Value = 500000000
$ws.Range("I51").Value = 500000000
$ws.Range("K51").Value = 500000000
$ws.Range("M51").Value = -499999516

$ws.Range("H53").Value = 544.2105
$ws.Range("J53").Value = 595.5
$ws.Range("L53").Value = 595.5
$ws.Range("N53").Value = -1869.5

$ws.Range("H62").Value = 3426.2856
$ws.Range("I62").Value = 2799.6
$ws.Range("J62").Value = 4993
$ws.Range("K62").Value = 2799.6
$ws.Range("L62").Value = 4993
$ws.Range("M62").Value = -2175.6
$ws.Range("N62").Value = -6241

$ws.Range("H65").Value = 3426.2856
$ws.Range("I65").Value = 2799.6
$ws.Range("J65").Value = 4993
$ws.Range("K65").Value = 13998
$ws.Range("L65").Value = 24965
$ws.Range("M65").Value = -10878
$ws.Range("N65").Value = -31205

$ws.Range("H98").Value = 1874.138
$ws.Range("I98").Value = 1933.9286
$ws.Range("K98").Value = 1933.9286
$ws.Range("M98").Value = -435.9286

$ws.Range("H116").Value = 8334.615
$ws.Range("J116").Value = 12182.667
$ws.Range("L116").Value = 12182.667
$ws.Range("N116").Value = -19066.667

$ws.Range("H122").Value = 1874.138
$ws.Range("I122").Value = 1933.9286
$ws.Range("K122").Value = 5801.7858
$ws.Range("M122").Value = -3351.7858

$ws.Range("H132").Value = 3345.647
$ws.Range("I132").Value = 3777.1592
$ws.Range("K132").Value = 11331.4776
$ws.Range("M132").Value = -8801.4776

$ws.Range("H135").Value = 913.8205
$ws.Range("I135").Value = 735.03125
$ws.Range("K135").Value = 6615.28125
$ws.Range("M135").Value = -4080.28125

$ws.Range("H137").Value = 1927579.6
$ws.Range("I137").Value = 2503705.5
$ws.Range("J137").Value = 7160.3335
$ws.Range("K137").Value = 7511116.5
$ws.Range("L137").Value = 21481.0005
$ws.Range("M137").Value = -7508566.5
$ws.Range("N137").Value = -26581.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 37075.082
$ws.Range("I45").Value = 48623.445
$ws.Range("J45").Value = 2430
$ws.Range("K45").Value = 48623.445
$ws.Range("L45").Value = 2430
$ws.Range("M45").Value = -48246.445
$ws.Range("N45").Value = -3184

$ws.Range("H61").Value = 2656.946
$ws.Range("I61").Value = 1128.3
$ws.Range("K61").Value = 1128.3
$ws.Range("M61").Value = -916.3

$ws.Range("H63").Value = 3814.6667
$ws.Range("I63").Value = 1630
$ws.Range("K63").Value = 1630
$ws.Range("M63").Value = -944

$ws.Range("H66").Value = 3814.6667
$ws.Range("I66").Value = 1630
$ws.Range("K66").Value = 8150
$ws.Range("M66").Value = -4718

$ws.Range("H74").Value = 193322.06
$ws.Range("I74").Value = 232930.83
$ws.Range("K74").Value = 232930.83
$ws.Range("M74").Value = -232056.83

$ws.Range("H77").Value = 193322.06
$ws.Range("I77").Value = 232930.83
$ws.Range("K77").Value = 1164654.15
$ws.Range("M77").Value = -1160286.15

$ws.Range("H102").Value = 2095.7058
$ws.Range("I102").Value = 1821.5
$ws.Range("K102").Value = 1821.5
$ws.Range("M102").Value = -199.5

$ws.Range("H122").Value = 4710.0625
$ws.Range("I122").Value = 4935.5186
$ws.Range("J122").Value = 3492.6
$ws.Range("K122").Value = 14806.5558
$ws.Range("L122").Value = 10477.8
$ws.Range("M122").Value = -12356.5558
$ws.Range("N122").Value = -15377.8

$ws.Range("H132").Value = 3104.7817
$ws.Range("I132").Value = 1993.8108
$ws.Range("K132").Value = 5981.4324
$ws.Range("M132").Value = -3451.4324

$ws.Range("H136").Value = 2656.946
$ws.Range("I136").Value = 1128.3
$ws.Range("K136").Value = 3384.9
$ws.Range("M136").Value = -834.8999999999996

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 20837870
$ws.Range("I20").Value = 29416710
$ws.Range("K20").Value = 29416710
$ws.Range("M20").Value = -29416463

$ws.Range("H94").Value = 42554410
$ws.Range("I94").Value = 64517212
$ws.Range("J94").Value = 1470
$ws.Range("K94").Value = 64517212
$ws.Range("L94").Value = 1470
$ws.Range("M94").Value = -64516761
$ws.Range("N94").Value = -2372

$ws.Range("H99").Value = 2929.3333
$ws.Range("I99").Value = 2265.077
$ws.Range("K99").Value = 2265.077
$ws.Range("M99").Value = -767.0770000000002

$ws.Range("H105").Value = 7648441.5
$ws.Range("I105").Value = 527460.4
$ws.Range("J105").Value = 16668351
$ws.Range("K105").Value = 527460.4
$ws.Range("L105").Value = 16668351
$ws.Range("M105").Value = -525713.4
$ws.Range("N105").Value = -16671845

$ws.Range("H107").Value = 1877421.8
$ws.Range("I107").Value = 2748480.5
$ws.Range("J107").Value = 1294.8462
$ws.Range("K107").Value = 2748480.5
$ws.Range("L107").Value = 1294.8462
$ws.Range("M107").Value = -2746560.5
$ws.Range("N107").Value = -5134.8462

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1509.9
$ws.Range("I22").Value = 1443.625
$ws.Range("J22").Value = 1775
$ws.Range("K22").Value = 1443.625
$ws.Range("L22").Value = 1775
$ws.Range("M22").Value = -1093.625
$ws.Range("N22").Value = -2475

$ws.Range("H31").Value = 2319912.8
$ws.Range("I31").Value = 3918.0625
$ws.Range("K31").Value = 3918.0625
$ws.Range("M31").Value = -3623.0625

$ws.Range("H34").Value = 2319912.8
$ws.Range("I34").Value = 3918.0625
$ws.Range("K34").Value = 3918.0625
$ws.Range("M34").Value = -3716.0625

$ws.Range("H62").Value = 20004394
$ws.Range("J62").Value = 3969
$ws.Range("L62").Value = 3969
$ws.Range("N62").Value = -5217

$ws.Range("H65").Value = 20004394
$ws.Range("J65").Value = 3969
$ws.Range("L65").Value = 19845
$ws.Range("N65").Value = -26085

$ws.Range("H105").Value = 742.9091
$ws.Range("I105").Value = 777.2
$ws.Range("K105").Value = 777.2
$ws.Range("M105").Value = 969.8

$ws.Range("H132").Value = 9261811
$ws.Range("I132").Value = 1240.9615
$ws.Range("K132").Value = 3722.8845
$ws.Range("M132").Value = -1192.8845

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 448.33334
$ws.Range("I107").Value = 637
$ws.Range("J107").Value = 354
$ws.Range("K107").Value = 1911
$ws.Range("L107").Value = 1062
$ws.Range("M107").Value = 9
$ws.Range("N107").Value = -4902

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2849.3076
$ws.Range("I122").Value = 2648.4688
$ws.Range("K122").Value = 7945.4064
$ws.Range("M122").Value = -5495.4064

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3847.2258
$ws.Range("I7").Value = 3224.1875
$ws.Range("J7").Value = 4511.8
$ws.Range("K7").Value = 3224.1875
$ws.Range("L7").Value = 4511.8
$ws.Range("M7").Value = -3112.1875
$ws.Range("N7").Value = -4735.8

$ws.Range("H22").Value = 2741.3333
$ws.Range("I22").Value = 2862.5
$ws.Range("K22").Value = 2862.5
$ws.Range("M22").Value = -2567.5

$ws.Range("H27").Value = 2741.3333
$ws.Range("I27").Value = 2862.5
$ws.Range("K27").Value = 2862.5
$ws.Range("M27").Value = -2755.5

$ws.Range("H68").Value = 6349.5
$ws.Range("I68").Value = 5133
$ws.Range("J68").Value = 9999
$ws.Range("K68").Value = 5133
$ws.Range("L68").Value = 9999
$ws.Range("M68").Value = -4384
$ws.Range("N68").Value = -11497

$ws.Range("H71").Value = 6349.5
$ws.Range("I71").Value = 5133
$ws.Range("J71").Value = 9999
$ws.Range("K71").Value = 25665
$ws.Range("L71").Value = 49995
$ws.Range("M71").Value = -21921
$ws.Range("N71").Value = -57483

$ws.Range("H122").Value = 6725.095
$ws.Range("I122").Value = 4895.8823
$ws.Range("J122").Value = 14499.25
$ws.Range("K122").Value = 14687.6469
$ws.Range("L122").Value = 43497.75
$ws.Range("M122").Value = -12237.6469
$ws.Range("N122").Value = -48397.75

$ws.Range("H126").Value = 3847.2258
$ws.Range("I126").Value = 3224.1875
$ws.Range("J126").Value = 4511.8
$ws.Range("K126").Value = 9672.5625
$ws.Range("L126").Value = 13535.4
$ws.Range("M126").Value = -7202.5625
$ws.Range("N126").Value = -18475.4

$ws.Range("H132").Value = 6114.1904
$ws.Range("I132").Value = 3338.0625
$ws.Range("K132").Value = 10014.1875
$ws.Range("M132").Value = -7484.1875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9321.333000000001
$ws.Range("J62").Value = 9984.571
$ws.Range("L62").Value = 9984.571
$ws.Range("N62").Value = -11232.571

$ws.Range("H65").Value = 9321.333000000001
$ws.Range("J65").Value = 9984.571
$ws.Range("L65").Value = 49922.855
$ws.Range("N65").Value = -56162.855

$ws.Range("H100").Value = 37037564
$ws.Range("I100").Value = 453.57895
$ws.Range("K100").Value = 907.1579
$ws.Range("M100").Value = -366.1579

$ws.Range("H122").Value = 10418127
$ws.Range("I122").Value = 1582.7368
$ws.Range("J122").Value = 50000990
$ws.Range("K122").Value = 4748.2104
$ws.Range("L122").Value = 150002970
$ws.Range("M122").Value = -2298.2104
$ws.Range("N122").Value = -150007870

$ws.Range("H126").Value = 13022.444
$ws.Range("I126").Value = 18533.834
$ws.Range("K126").Value = 55601.50199999999
$ws.Range("M126").Value = -53131.50199999999

$ws.Range("H136").Value = 11277.842
$ws.Range("I136").Value = 11208.229
$ws.Range("K136").Value = 33624.687
$ws.Range("M136").Value = -31074.687
